$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tb_rkap")

# Insert a new column D ("uraian") between gl_acc (C) and tahun_rkap (old D, now E).
$ws.Columns("D").Insert()

# Row values (tb_rkap_detail descriptions) - filled in before the header,
# matching the order the values were looked up / entered in the source workbook
$ws.Range("D2").Value = "Suku Cadang"
$ws.Range("D3").Value = "Bahan Pembantu"
$ws.Range("D5").Value = "Pemeliharaan - Kendaraan"
$ws.Range("D6").Value = "Pajak Kendaraan"
$ws.Range("D7").Value = "Retribusi - Kendaraan"
$ws.Range("D9").Value = "Sewa - Kendaraan"
$ws.Range("D4").Value = "Lainnya - Bahan Bakar Kendaraan"
$ws.Range("D8").Value = "Retribusi - Parkir & Tol"

# Header
$ws.Range("D1").Value = "uraian"

# Resize the new column to fit its contents.
$ws.Columns("D").AutoFit()

# Move the selection to the newly added column, matching the saved view state.
$ws.Range("D2").Select()
